$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 11.071066
$ws.Cells.Item(2, 8).Value = 33.213198
$ws.Cells.Item(2, 9).Value = 0.07439893415073216
$ws.Cells.Item(2, 10).Value = 0.07645010820826473
$ws.Cells.Item(2, 13).Value = 2.266198333333334
$ws.Cells.Item(2, 14).Value = 6.798595000000001
$ws.Cells.Item(2, 15).Value = 0.004051889747168086
$ws.Cells.Item(2, 16).Value = 0.004065587792570741
$ws.Cells.Item(2, 17).Value = 25.08923131742334
$ws.Cells.Item(2, 18).Value = 225.80308185681
$ws.Cells.Item(2, 19).Value = 0.0003014562784855853
$ws.Cells.Item(2, 20).Value = 0.0003108146266722333
$ws.Cells.Item(3, 7).Value = 11.071066
$ws.Cells.Item(3, 8).Value = 33.213198
$ws.Cells.Item(3, 9).Value = 0.07439893415073216
$ws.Cells.Item(3, 10).Value = 0.07645010820826473
$ws.Cells.Item(3, 15).Value = 0.00469498661713519
$ws.Cells.Item(3, 16).Value = 0.004710858751832666
$ws.Cells.Item(3, 17).Value = 29.07127602665867
$ws.Cells.Item(3, 18).Value = 261.641484239928
$ws.Cells.Item(3, 19).Value = 0.0003493020001668098
$ws.Cells.Item(3, 20).Value = 0.0003601456613314582
$ws.Cells.Item(4, 7).Value = 11.071066
$ws.Cells.Item(4, 8).Value = 33.213198
$ws.Cells.Item(4, 9).Value = 0.07439893415073216
$ws.Cells.Item(4, 10).Value = 0.07645010820826473
$ws.Cells.Item(4, 13).Value = 258.491684
$ws.Cells.Item(4, 14).Value = 775.475052
$ws.Cells.Item(4, 15).Value = 0.462174818824101
$ws.Cells.Item(4, 16).Value = 0.4637372729004096
$ws.Cells.Item(4, 17).Value = 2861.778494015144
$ws.Cells.Item(4, 18).Value = 25756.0064461363
$ws.Cells.Item(4, 19).Value = 0.03438531391182086
$ws.Cells.Item(4, 20).Value = 0.0354527646934419
$ws.Cells.Item(5, 7).Value = 11.071066
$ws.Cells.Item(5, 8).Value = 33.213198
$ws.Cells.Item(5, 9).Value = 0.07439893415073216
$ws.Cells.Item(5, 10).Value = 0.07645010820826473
$ws.Cells.Item(5, 13).Value = 5.653232
$ws.Cells.Item(5, 14).Value = 11.306464
$ws.Cells.Item(5, 15).Value = 0.01010779702828123
$ws.Cells.Item(5, 16).Value = 0.006761312008663634
$ws.Cells.Item(5, 17).Value = 62.587304585312
$ws.Cells.Item(5, 18).Value = 375.523827511872
$ws.Cells.Item(5, 19).Value = 0.0007520093255160611
$ws.Cells.Item(5, 20).Value = 0.0005169030346921746
$ws.Cells.Item(6, 7).Value = 11.071066
$ws.Cells.Item(6, 8).Value = 33.213198
$ws.Cells.Item(6, 9).Value = 0.07439893415073216
$ws.Cells.Item(6, 10).Value = 0.07645010820826473
$ws.Cells.Item(6, 13).Value = 290.257182
$ws.Cells.Item(6, 14).Value = 870.7715460000001
$ws.Cells.Item(6, 15).Value = 0.5189705077833145
$ws.Cells.Item(6, 16).Value = 0.5207249685465234
$ws.Cells.Item(6, 17).Value = 3213.456418896012
$ws.Cells.Item(6, 18).Value = 28921.10777006411
$ws.Cells.Item(6, 19).Value = 0.03861085263474285
$ws.Cells.Item(6, 20).Value = 0.03980948019212696
$ws.Cells.Item(7, 9).Value = 0.7926488104270929
$ws.Cells.Item(7, 10).Value = 0.8145020895800995
$ws.Cells.Item(7, 13).Value = 2.266198333333334
$ws.Cells.Item(7, 14).Value = 6.798595000000001
$ws.Cells.Item(7, 15).Value = 0.004051889747168086
$ws.Cells.Item(7, 16).Value = 0.004065587792570741
$ws.Cells.Item(7, 17).Value = 267.3015357719351
$ws.Cells.Item(7, 18).Value = 2405.713821947415
$ws.Cells.Item(7, 19).Value = 0.003211725588074518
$ws.Cells.Item(7, 20).Value = 0.003311429752420212
$ws.Cells.Item(8, 9).Value = 0.7926488104270929
$ws.Cells.Item(8, 10).Value = 0.8145020895800995
$ws.Cells.Item(8, 15).Value = 0.00469498661713519
$ws.Cells.Item(8, 16).Value = 0.004710858751832666
$ws.Cells.Item(8, 19).Value = 0.00372147555704333
$ws.Cells.Item(8, 20).Value = 0.003837004297084405
$ws.Cells.Item(9, 9).Value = 0.7926488104270929
$ws.Cells.Item(9, 10).Value = 0.8145020895800995
$ws.Cells.Item(9, 13).Value = 258.491684
$ws.Cells.Item(9, 14).Value = 775.475052
$ws.Cells.Item(9, 15).Value = 0.462174818824101
$ws.Cells.Item(9, 16).Value = 0.4637372729004096
$ws.Cells.Item(9, 17).Value = 30489.486776668
$ws.Cells.Item(9, 18).Value = 274405.380990012
$ws.Cells.Item(9, 19).Value = 0.3663423203502809
$ws.Cells.Item(9, 20).Value = 0.3777149777935604
$ws.Cells.Item(10, 9).Value = 0.7926488104270929
$ws.Cells.Item(10, 10).Value = 0.8145020895800995
$ws.Cells.Item(10, 13).Value = 5.653232
$ws.Cells.Item(10, 14).Value = 11.306464
$ws.Cells.Item(10, 15).Value = 0.01010779702828123
$ws.Cells.Item(10, 16).Value = 0.006761312008663634
$ws.Cells.Item(10, 17).Value = 666.807301659408
$ws.Cells.Item(10, 18).Value = 4000.843809956448
$ws.Cells.Item(10, 19).Value = 0.008011933290505618
$ws.Cells.Item(10, 20).Value = 0.00550710275935955
$ws.Cells.Item(11, 9).Value = 0.7926488104270929
$ws.Cells.Item(11, 10).Value = 0.8145020895800995
$ws.Cells.Item(11, 13).Value = 290.257182
$ws.Cells.Item(11, 14).Value = 870.7715460000001
$ws.Cells.Item(11, 15).Value = 0.5189705077833145
$ws.Cells.Item(11, 16).Value = 0.5207249685465234
$ws.Cells.Item(11, 17).Value = 34236.27551755946
$ws.Cells.Item(11, 18).Value = 308126.4796580352
$ws.Cells.Item(11, 19).Value = 0.4113613556411886
$ws.Cells.Item(11, 20).Value = 0.4241315749776749
$ws.Cells.Item(12, 7).Value = 5.650614
$ws.Cells.Item(12, 8).Value = 16.951842
$ws.Cells.Item(12, 9).Value = 0.03797282564273442
$ws.Cells.Item(12, 10).Value = 0.03901973411983414
$ws.Cells.Item(12, 13).Value = 2.266198333333334
$ws.Cells.Item(12, 14).Value = 6.798595000000001
$ws.Cells.Item(12, 15).Value = 0.004051889747168086
$ws.Cells.Item(12, 16).Value = 0.004065587792570741
$ws.Cells.Item(12, 17).Value = 12.80541202911
$ws.Cells.Item(12, 18).Value = 115.24870826199
$ws.Cells.Item(12, 19).Value = 0.000153861702892797
$ws.Cells.Item(12, 20).Value = 0.0001586381547069537
$ws.Cells.Item(13, 7).Value = 5.650614
$ws.Cells.Item(13, 8).Value = 16.951842
$ws.Cells.Item(13, 9).Value = 0.03797282564273442
$ws.Cells.Item(13, 10).Value = 0.03901973411983414
$ws.Cells.Item(13, 15).Value = 0.00469498661713519
$ws.Cells.Item(13, 16).Value = 0.004710858751832666
$ws.Cells.Item(13, 17).Value = 14.837826756168
$ws.Cells.Item(13, 18).Value = 133.540440805512
$ws.Cells.Item(13, 19).Value = 0.0001782819082074461
$ws.Cells.Item(13, 20).Value = 0.0001838164559726043
$ws.Cells.Item(14, 7).Value = 5.650614
$ws.Cells.Item(14, 8).Value = 16.951842
$ws.Cells.Item(14, 9).Value = 0.03797282564273442
$ws.Cells.Item(14, 10).Value = 0.03901973411983414
$ws.Cells.Item(14, 13).Value = 258.491684
$ws.Cells.Item(14, 14).Value = 775.475052
$ws.Cells.Item(14, 15).Value = 0.462174818824101
$ws.Cells.Item(14, 16).Value = 0.4637372729004096
$ws.Cells.Item(14, 17).Value = 1460.636728493976
$ws.Cells.Item(14, 18).Value = 13145.73055644578
$ws.Cells.Item(14, 19).Value = 0.01755008381166996
$ws.Cells.Item(14, 20).Value = 0.01809490509003095
$ws.Cells.Item(15, 7).Value = 5.650614
$ws.Cells.Item(15, 8).Value = 16.951842
$ws.Cells.Item(15, 9).Value = 0.03797282564273442
$ws.Cells.Item(15, 10).Value = 0.03901973411983414
$ws.Cells.Item(15, 13).Value = 5.653232
$ws.Cells.Item(15, 14).Value = 11.306464
$ws.Cells.Item(15, 15).Value = 0.01010779702828123
$ws.Cells.Item(15, 16).Value = 0.006761312008663634
$ws.Cells.Item(15, 17).Value = 31.944231884448
$ws.Cells.Item(15, 18).Value = 191.665391306688
$ws.Cells.Item(15, 19).Value = 0.0003838216141870721
$ws.Cells.Item(15, 20).Value = 0.0002638245968792967
$ws.Cells.Item(16, 7).Value = 5.650614
$ws.Cells.Item(16, 8).Value = 16.951842
$ws.Cells.Item(16, 9).Value = 0.03797282564273442
$ws.Cells.Item(16, 10).Value = 0.03901973411983414
$ws.Cells.Item(16, 13).Value = 290.257182
$ws.Cells.Item(16, 14).Value = 870.7715460000001
$ws.Cells.Item(16, 15).Value = 0.5189705077833145
$ws.Cells.Item(16, 16).Value = 0.5207249685465234
$ws.Cells.Item(16, 17).Value = 1640.131296209748
$ws.Cells.Item(16, 18).Value = 14761.18166588773
$ws.Cells.Item(16, 19).Value = 0.01970677660577715
$ws.Cells.Item(16, 20).Value = 0.02031854982224434
$ws.Cells.Item(17, 7).Value = 11.977561
$ws.Cells.Item(17, 8).Value = 23.955122
$ws.Cells.Item(17, 9).Value = 0.08049069277749564
$ws.Cells.Item(17, 10).Value = 0.05513987749816152
$ws.Cells.Item(17, 13).Value = 2.266198333333334
$ws.Cells.Item(17, 14).Value = 6.798595000000001
$ws.Cells.Item(17, 15).Value = 0.004051889747168086
$ws.Cells.Item(17, 16).Value = 0.004065587792570741
$ws.Cells.Item(17, 17).Value = 27.14352877559833
$ws.Cells.Item(17, 18).Value = 162.86117265359
$ws.Cells.Item(17, 19).Value = 0.0003261394128075909
$ws.Cells.Item(17, 20).Value = 0.0002241760128403715
$ws.Cells.Item(18, 7).Value = 11.977561
$ws.Cells.Item(18, 8).Value = 23.955122
$ws.Cells.Item(18, 9).Value = 0.08049069277749564
$ws.Cells.Item(18, 10).Value = 0.05513987749816152
$ws.Cells.Item(18, 15).Value = 0.00469498661713519
$ws.Cells.Item(18, 16).Value = 0.004710858751832666
$ws.Cells.Item(18, 17).Value = 31.45162190859867
$ws.Cells.Item(18, 18).Value = 188.709731451592
$ws.Cells.Item(18, 19).Value = 0.0003779027253942821
$ws.Cells.Item(18, 20).Value = 0.0002597561744871952
$ws.Cells.Item(19, 7).Value = 11.977561
$ws.Cells.Item(19, 8).Value = 23.955122
$ws.Cells.Item(19, 9).Value = 0.08049069277749564
$ws.Cells.Item(19, 10).Value = 0.05513987749816152
$ws.Cells.Item(19, 13).Value = 258.491684
$ws.Cells.Item(19, 14).Value = 775.475052
$ws.Cells.Item(19, 15).Value = 0.462174818824101
$ws.Cells.Item(19, 16).Value = 0.4637372729004096
$ws.Cells.Item(19, 17).Value = 3096.099913102724
$ws.Cells.Item(19, 18).Value = 18576.59947861634
$ws.Cells.Item(19, 19).Value = 0.03720077135146542
$ws.Cells.Item(19, 20).Value = 0.02557041641906008
$ws.Cells.Item(20, 7).Value = 11.977561
$ws.Cells.Item(20, 8).Value = 23.955122
$ws.Cells.Item(20, 9).Value = 0.08049069277749564
$ws.Cells.Item(20, 10).Value = 0.05513987749816152
$ws.Cells.Item(20, 13).Value = 5.653232
$ws.Cells.Item(20, 14).Value = 11.306464
$ws.Cells.Item(20, 15).Value = 0.01010779702828123
$ws.Cells.Item(20, 16).Value = 0.006761312008663634
$ws.Cells.Item(20, 17).Value = 67.711931127152
$ws.Cells.Item(20, 18).Value = 270.847724508608
$ws.Cells.Item(20, 19).Value = 0.0008135835852606676
$ws.Cells.Item(20, 20).Value = 0.0003728179158845612
$ws.Cells.Item(21, 7).Value = 11.977561
$ws.Cells.Item(21, 8).Value = 23.955122
$ws.Cells.Item(21, 9).Value = 0.08049069277749564
$ws.Cells.Item(21, 10).Value = 0.05513987749816152
$ws.Cells.Item(21, 13).Value = 290.257182
$ws.Cells.Item(21, 14).Value = 870.7715460000001
$ws.Cells.Item(21, 15).Value = 0.5189705077833145
$ws.Cells.Item(21, 16).Value = 0.5207249685465234
$ws.Cells.Item(21, 17).Value = 3476.573103093102
$ws.Cells.Item(21, 18).Value = 20859.43861855861
$ws.Cells.Item(21, 19).Value = 0.04177229570256768
$ws.Cells.Item(21, 20).Value = 0.02871271097588931
$ws.Cells.Item(22, 7).Value = 2.156022333333333
$ws.Cells.Item(22, 8).Value = 6.468067
$ws.Cells.Item(22, 9).Value = 0.01448873700194494
$ws.Cells.Item(22, 10).Value = 0.0148881905936401
$ws.Cells.Item(22, 13).Value = 2.266198333333334
$ws.Cells.Item(22, 14).Value = 6.798595000000001
$ws.Cells.Item(22, 15).Value = 0.004051889747168086
$ws.Cells.Item(22, 16).Value = 0.004065587792570741
$ws.Cells.Item(22, 17).Value = 4.885974218429444
$ws.Cells.Item(22, 18).Value = 43.973767965865
$ws.Cells.Item(22, 19).Value = 0.00005870676490759558
$ws.Cells.Item(22, 20).Value = 0.00006052924593096972
$ws.Cells.Item(23, 7).Value = 2.156022333333333
$ws.Cells.Item(23, 8).Value = 6.468067
$ws.Cells.Item(23, 9).Value = 0.01448873700194494
$ws.Cells.Item(23, 10).Value = 0.0148881905936401
$ws.Cells.Item(23, 15).Value = 0.00469498661713519
$ws.Cells.Item(23, 16).Value = 0.004710858751832666
$ws.Cells.Item(23, 17).Value = 5.661453049956888
$ws.Cells.Item(23, 18).Value = 50.953077449612
$ws.Cells.Item(23, 19).Value = 0.00006802442632332293
$ws.Cells.Item(23, 20).Value = 0.00007013616295700224
$ws.Cells.Item(24, 7).Value = 2.156022333333333
$ws.Cells.Item(24, 8).Value = 6.468067
$ws.Cells.Item(24, 9).Value = 0.01448873700194494
$ws.Cells.Item(24, 10).Value = 0.0148881905936401
$ws.Cells.Item(24, 13).Value = 258.491684
$ws.Cells.Item(24, 14).Value = 775.475052
$ws.Cells.Item(24, 15).Value = 0.462174818824101
$ws.Cells.Item(24, 16).Value = 0.4637372729004096
$ws.Cells.Item(24, 17).Value = 557.3138436849426
$ws.Cells.Item(24, 18).Value = 5015.824593164483
$ws.Cells.Item(24, 19).Value = 0.006696329398863951
$ws.Cells.Item(24, 20).Value = 0.006904208904316191
$ws.Cells.Item(25, 7).Value = 2.156022333333333
$ws.Cells.Item(25, 8).Value = 6.468067
$ws.Cells.Item(25, 9).Value = 0.01448873700194494
$ws.Cells.Item(25, 10).Value = 0.0148881905936401
$ws.Cells.Item(25, 13).Value = 5.653232
$ws.Cells.Item(25, 14).Value = 11.306464
$ws.Cells.Item(25, 15).Value = 0.01010779702828123
$ws.Cells.Item(25, 16).Value = 0.006761312008663634
$ws.Cells.Item(25, 17).Value = 12.18849444751467
$ws.Cells.Item(25, 18).Value = 73.130966685088
$ws.Cells.Item(25, 19).Value = 0.0001464492128118073
$ws.Cells.Item(25, 20).Value = 0.0001006637018480518
$ws.Cells.Item(26, 7).Value = 2.156022333333333
$ws.Cells.Item(26, 8).Value = 6.468067
$ws.Cells.Item(26, 9).Value = 0.01448873700194494
$ws.Cells.Item(26, 10).Value = 0.0148881905936401
$ws.Cells.Item(26, 13).Value = 290.257182
$ws.Cells.Item(26, 14).Value = 870.7715460000001
$ws.Cells.Item(26, 15).Value = 0.5189705077833145
$ws.Cells.Item(26, 16).Value = 0.5207249685465234
$ws.Cells.Item(26, 17).Value = 3213.456418896012
$ws.Cells.Item(26, 18).Value = 14761.18166588773
$ws.Cells.Item(26, 19).Value = 0.007519227199038263
$ws.Cells.Item(26, 20).Value = 0.007752652578587887
